# Fixed blank node bug with jena pretty format, adapted formatting
# Add a new time-tracking entry (row 50) to the "Tabelle1" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New row: Datum / Stunden / Aufgabe
$ws.Range("A50").Value = 45356
$ws.Range("A50").NumberFormat = $ws.Range("A49").NumberFormat
$ws.Range("B50").Value = 3
$ws.Range("C50").Value = "Bug Blank Nodes Anzeige"

# Update the active selection to follow the newly added row, like Excel would.
$ws.Range("C50").Select() | Out-Null

Write-Host "Added row 50 to Tabelle1"
